# Update the Notes sheet: correct the wording of the existing biomark-antenna
# note and add a new note about frontend/backend antenna-site-code entry counts.
# This mirrors the author editing text directly in the "Notes" worksheet
# (the 4th sheet in the workbook) and leaving it as the active sheet on save.

$wb = $excel.ActiveWorkbook
$notes = $wb.Worksheets.Item("Notes")

$notes.Range("A2").Value2 = "For the biomark antennas, I don't think you can actually code the readers to detect as the frontend codes; only codes like A4,B2 etc…"
$notes.Range("A3").Value2 = "For antenna metadatea codes, the frontend code should only have 1 entry. It's ok for the biomark backend codes to have multiple entries, but not the stationary ones"

# Make "Notes" the active sheet / active tab, with the selection resting on A1.
$notes.Activate()
$notes.Cells.Item(1, 1).Select()
